# update search and upload
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Expected Login (column E) for rows 2, 4, 5: "Blank" -> guest-access message
$ws.Range("E2").Value = "You are currently using guest access\nLog in"
$ws.Range("E4").Value = "You are currently using guest access\nLog in"
$ws.Range("E5").Value = "You are currently using guest access\nLog in"

# Update the selected cell shown in the sheet view
$ws.Range("D14").Select()
